# Updated symbol list on Tue Dec 20 23:58:10 UTC 2022 with GitHub Actions
#
# This script applies the updated crypto price/volume figures to Sheet1.
# Numeric-looking values are stored in this workbook as plain text, so we
# force a Text number format before assigning them to avoid Excel silently
# re-typing the cell as a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$Address, [string]$Value)
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
}

# --- simple price updates (column D) ---
Set-TextValue $ws "D2"  "251.78"
Set-TextValue $ws "D4"  "5.438"
Set-TextValue $ws "D5"  "0.05711"
Set-TextValue $ws "D6"  "3.428"
Set-TextValue $ws "D7"  "6.367"
Set-TextValue $ws "D8"  "0.8128"
Set-TextValue $ws "D9"  "0.9366"
Set-TextValue $ws "D11" "0.07431"
Set-TextValue $ws "D13" "0.03077"
Set-TextValue $ws "D14" "0.09367"
Set-TextValue $ws "D15" "3.715"
Set-TextValue $ws "D16" "0.001586"
Set-TextValue $ws "D17" "0.04777"

# column E text update (row 18)
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextValue $ws "D20" "0.005041"
Set-TextValue $ws "D21" "0.001026"
Set-TextValue $ws "D22" "0.0001499"
Set-TextValue $ws "D23" "3.705"
Set-TextValue $ws "D24" "2.170"

Set-TextValue $ws "D40" "0.04028"

# --- rows 41-43 got reordered/rewritten (Kick/BKEX/CEJI shuffled) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006788"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1072"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002709"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue $ws "D45" "0.00005799"
